$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("FIRST3").Name = "FIRST"
$wb.Worksheets.Item("SECOND3").Name = "SECOND"
$wb.Worksheets.Item("THIRD3").Name = "THIRD"
